{"js": "// Rename the \"estadoApp\" field to \"aceptado\" and reword its description\n// (spans the two lines following the field name) for the Jugador\n// collection, per commit \"Create new field 'Aceotado' in Equipo and\n// Jugador collection\".\n\n// Locate the paragraph that still has the old field name \"estadoApp\" so\n// the script doesn't depend on a hard-coded paragraph index.\nconst hits = context.document.body.search(\"estadoApp\", { matchCase: true });\nhits.load(\"items\");\nawait context.sync();\n\nif (hits.items.length === 0) {\n  throw new Error('Could not find \"estadoApp\" field to rename.');\n}\n\nconst firstLine = hits.items[0].paragraphs.getFirst();\nconst secondLine = firstLine.getNext();\nfirstLine.load(\"text\");\nsecondLine.load(\"text\");\nawait context.sync();\n\n// First line keeps its leading tab, only the field name/description text\n// changes: \"estadoApp: (Muestra ... en la\" -> \"aceptado: (Indica ... a la \".\nfirstLine.getRange().insertText(\n  \"\\taceptado: (Indica si el jugador ha sido aceptado en la liga o est\u00e1 a la \",\n  \"Replace\"\n);\n\n// Second line loses its leading tabs/spaces, is indented as a hanging\n// continuation instead, and the wording + closing parenthesis changes:\n// \"liga o est\u00e1 a la espera de serlo. Aceptado o pendiente) \"\n// -> \"      espera de ello) \" (indented via paragraph format).\nsecondLine.leftIndent = 70.8; // points == 1416 twips\nsecondLine.getRange().insertText(\"      espera de ello) \", \"Replace\");\n\nawait context.sync();\n", "ps1": "# Rename the \"estadoApp\" field to \"aceptado\" and reword its description\n# (spans the two lines following the field name) for the Jugador\n# collection, per commit \"Create new field 'Aceotado' in Equipo and\n# Jugador collection\".\n\n$d = $word.ActiveDocument\n\n# Locate the field by searching for its old name so the script doesn't\n# depend on a hard-coded paragraph index.\n$findRange = $d.Content\n$found = $findRange.Find.Execute(\"estadoApp\")\nif (-not $found) {\n    throw 'Could not find \"estadoApp\" field to rename.'\n}\n$targetStart = $findRange.Start\n\n$paragraphs = $d.Paragraphs\n$firstLine = $null\nfor ($i = 1; $i -le $paragraphs.Count; $i++) {\n    $candidate = $paragraphs.Item($i)\n    if ($candidate.Range.Start -le $targetStart -and $candidate.Range.End -ge $targetStart) {\n        $firstLine = $candidate\n        break\n    }\n}\nif ($null -eq $firstLine) {\n    throw \"Could not locate the paragraph containing the estadoApp field.\"\n}\n$secondLine = $firstLine.Next()\n\n# First line keeps its leading tab, only the field name/description text\n# changes: \"estadoApp: (Muestra ... en la\" -> \"aceptado: (Indica ... a la \".\n$r1 = $firstLine.Range\n$r1Body = $d.Range($r1.Start, $r1.End - 1)\n$r1Body.Text = \"`taceptado: (Indica si el jugador ha sido aceptado en la liga o est\u00e1 a la \"\n\n# Second line loses its leading tabs/spaces, is indented as a hanging\n# continuation instead, and the wording + closing parenthesis changes:\n# \"liga o est\u00e1 a la espera de serlo. Aceptado o pendiente) \"\n# -> \"      espera de ello) \" (indented via paragraph format).\n$secondLine.Format.LeftIndent = 70.8  # points == 1416 twips\n$r2 = $secondLine.Range\n$r2Body = $d.Range($r2.Start, $r2.End - 1)\n$r2Body.Text = \"      espera de ello) \"\n"}
